# [Feat] => Descarga formato pdf Lista de Trabajo en Busqueda Masiva
#
# The "Expedientes" report sheet is being repurposed for the new
# "Busqueda Masiva" (mass search) / Solicitudes flow, so it is renamed to
# "Solicitudes". Excel automatically re-points the existing defined names
# (Expedientes / Expedientes_Estudios) to the renamed sheet.

$wb = $excel.ActiveWorkbook

# Grab the sheet by its current name rather than assuming ActiveSheet,
# so the script is resilient regardless of which sheet/window is active
# when it runs.
$ws = $wb.Worksheets.Item("Expedientes")

$ws.Name = "Solicitudes"

# Leave the cursor where the author last left it when saving the file.
$ws.Range("E32").Select() | Out-Null
